# Insert a new weekly record at row 191 for "Terminal La Palmera de La Serena -
# Albahaca", pushing the existing rows 191-232 down to 192-233 (dimension grows
# from A1:R232 to A1:R233).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 191..232 down by one row.
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new weekly observation.
$ws.Range("A191").Value = 8
$ws.Range("B191").Value = "Terminal La Palmera de La Serena"
$ws.Range("C191").Value = "Coquimbo"
$ws.Range("D191").Value = 45275
$ws.Range("E191").Value = 4
$ws.Range("F191").Value = 100112052
$ws.Range("G191").Value = "Albahaca"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 1000
$ws.Range("K191").Value = 3400
$ws.Range("L191").Value = 3500
$ws.Range("M191").Value = 3450
$ws.Range("N191").Value = "$/docena de matas"
$ws.Range("O191").Value = "Provincia del Elquí"
$ws.Range("P191").Value = 575
$ws.Range("Q191").Value = 6
$ws.Range("R191").Value = "Hortaliza"
